$d = $word.ActiveDocument

# Replace "In processing." (with trailing period) occurrences with "Completed."
$d.Content.Find.Execute("In processing.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Completed.", 2)

# Replace remaining "In processing" (no trailing period) occurrences with "Completed."
$d.Content.Find.Execute("In processing", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Completed.", 2)
